# Überarbeitung für Kurstag 5
# - Clear a handful of stray/duplicate Konzentration (C) readings.
# - Fix the Temperatur (D) / Temperaturatur °F (E) values for the NaCl
#   (rows 52-101) and KOH (rows 102-151) blocks so they match the H2SO4
#   block: 25 °C / 77.0 °F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray Konzentration (mg/L) values in column C, keeping the
# cell itself present but empty (no cached value, default style).
$rowsToClearC = @(10, 19, 39, 65, 78, 108, 128, 142)
foreach ($r in $rowsToClearC) {
    $ws.Cells.Item($r, 3).ClearContents()
    $ws.Cells.Item($r, 3).Style = "Normal"
}

# NaCl block (rows 52-101): Temperatur 30 -> 25, Temperaturatur °F 86.0 -> 77.0
for ($r = 52; $r -le 101; $r++) {
    $ws.Cells.Item($r, 4).Value = 25
    $ws.Cells.Item($r, 5).Value = 77.0
}

# KOH block (rows 102-151): Temperatur 28 -> 25, Temperaturatur °F 82.4 -> 77.0
for ($r = 102; $r -le 151; $r++) {
    $ws.Cells.Item($r, 4).Value = 25
    $ws.Cells.Item($r, 5).Value = 77.0
}
